$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Taxonsorteringsordning bump ---
$ws.Range("B2").Value = 79244

# --- Rows 3 and 4: content effectively swaps between the two records, plus
#     the Taxonsorteringsordning (col B) values are bumped by one.
#     Write the full target state for each cell directly rather than
#     attempting an in-place swap.

# Row 3 (becomes what was row 4's record, with B incremented to 79245)
$ws.Range("A3").Value = 130853761
$ws.Range("B3").Value = 79245
$ws.Range("E3").Value = 230405
$ws.Range("F3").Value = "Garnlav (ssp. sarmentosa)"
$ws.Range("G3").Value = "Alectoria sarmentosa subsp. sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("P3").Value = "Djupbäcken, Djupbäcken, Jmt"
$ws.Range("Q3").Value = 442771
$ws.Range("R3").Value = 7039709
$ws.Range("S3").Value = 20
$ws.Range("Z3").Value = "11:05"
$ws.Range("AB3").Value = "11:05"
$ws.Range("AC3").ClearContents()
$ws.Range("AW3").Value = "Maria Danvind"
$ws.Range("AX3").Value = "Maria Danvind"

# Row 4 (becomes what was row 3's record, with B incremented to 91805)
$ws.Range("A4").Value = 130861152
$ws.Range("B4").Value = 91805
$ws.Range("E4").Value = 1108
$ws.Range("F4").Value = "Harticka"
$ws.Range("G4").Value = "Pelloporus leporinus"
$ws.Range("H4").Value = "(Fr.) Krieglst."
$ws.Range("P4").Value = "Djupbäcken, Jmt"
$ws.Range("Q4").Value = 442868
$ws.Range("R4").Value = 7039767
$ws.Range("S4").Value = 10
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
$ws.Range("AC4").Value = "I stående levande gran med full längd."
$ws.Range("AW4").Value = "Kristian Zackrisson"
$ws.Range("AX4").Value = "Kristian Zackrisson"

# --- Rows 5 and 6: coordinates + comment swap between the two records,
#     plus both records' Taxonsorteringsordning bump by one.
$ws.Range("A5").Value = 130861158
$ws.Range("B5").Value = 79244
$ws.Range("Q5").Value = 442743
$ws.Range("R5").Value = 7039650
$ws.Range("AC5").Value = "På gran."

$ws.Range("A6").Value = 130861156
$ws.Range("B6").Value = 79244
$ws.Range("Q6").Value = 442897
$ws.Range("R6").Value = 7039676
$ws.Range("AC6").Value = "På död stående gran med full längd."

# --- Remaining rows: simple Taxonsorteringsordning (col B) bumps by one ---
$ws.Range("B7").Value = 91805
$ws.Range("B8").Value = 91829
$ws.Range("B9").Value = 79244
$ws.Range("B10").Value = 79244
$ws.Range("B12").Value = 91829
$ws.Range("B13").Value = 79244
$ws.Range("B14").Value = 91805
$ws.Range("B16").Value = 79244
$ws.Range("B17").Value = 79244
$ws.Range("B18").Value = 79244
$ws.Range("B19").Value = 79244
